# Weekly refresh of the "Hortaliza, Terminal La Palmera de La Serena - Berenjena"
# data sheet: a new daily price record is inserted at row 257 (pushing the
# existing rows 257-308 down to 258-309), growing the sheet from 308 to 309
# data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 257; Excel shifts rows
# 257:308 down to 258:309 and carries the row-above formatting (incl. the
# date-number-format style used by column D) onto the new row.
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A257").Value = 8
$ws.Range("B257").Value = "Terminal La Palmera de La Serena"
$ws.Range("C257").Value = "Coquimbo"
$ws.Range("D257").Value = 45211
$ws.Range("E257").Value = 4
$ws.Range("F257").Value = 100112001
$ws.Range("G257").Value = "Berenjena"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 400
$ws.Range("K257").Value = 9000
$ws.Range("L257").Value = 9500
$ws.Range("M257").Value = 9250
$ws.Range("N257").Value = "`$/caja 50 unidades"
$ws.Range("O257").Value = "Región de Arica y Parinacota"
$ws.Range("P257").Value = 185
$ws.Range("Q257").Value = 50
$ws.Range("R257").Value = "Hortaliza"
